# Weekly fruit/vegetable update: insert a new price record as row 721
# (Naranja - Valencia, Macroferia Regional de Talca), pushing the existing
# rows 721:756 down to 722:757.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 721; Excel shifts rows 721:756
# down to 722:757 and keeps the sheet's date-format style on column D.
$ws.Rows.Item(721).Insert()

$ws.Cells.Item(721, 1).Value = 5
$ws.Cells.Item(721, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(721, 3).Value = "Maule"
$ws.Cells.Item(721, 4).Value = 44939
$ws.Cells.Item(721, 5).Value = 7
$ws.Cells.Item(721, 6).Value = "Fruta"
$ws.Cells.Item(721, 7).Value = 100102
$ws.Cells.Item(721, 8).Value = "Cítricos"
$ws.Cells.Item(721, 9).Value = 100102005
$ws.Cells.Item(721, 10).Value = "Naranja"
$ws.Cells.Item(721, 11).Value = "Valencia"
$ws.Cells.Item(721, 12).Value = "Primera"
$ws.Cells.Item(721, 13).Value = 350
$ws.Cells.Item(721, 14).Value = 12000
$ws.Cells.Item(721, 15).Value = 12000
$ws.Cells.Item(721, 16).Value = 12000
$ws.Cells.Item(721, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(721, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(721, 19).Value = 800
$ws.Cells.Item(721, 20).Value = 15
